# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
# Rework the "Source:" block at the bottom of the Uganda summary sheet:
#  - drop the hyperlink on the UBOS report URL and turn it into a second
#    plain-text "Uganda Bureau of Statistics" line
#  - add the URL itself back in as a plain text line below it
#  - replace the two long citation paragraphs with short "UBS" labels

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that currently lives on A64 (and its relationship).
$ws.Range("A64").Hyperlinks.Delete()

# --- "Source:" block (rows 62-66) ---------------------------------------
# Rows 62/63/65 already use the italic "source" look and keep using it, so
# only their text needs touching (leaves existing formatting untouched).
$ws.Range("A62").Value = "Source:"
$ws.Range("A63").Value = ""
$ws.Range("A65").Value = ""

# Row 64 used to be the hyperlinked URL; it now just repeats the
# "Uganda Bureau of Statistics" text as plain (no link) italic text, so the
# HyperLink font (underline+blue) must be swapped for the italic "source" one.
$ws.Range("A64").Value = "Uganda Bureau of Statistics"
$ws.Range("A64").Font.Underline = $false
$ws.Range("A64").Font.Italic = $true

# Row 66 is new: the bare report URL as plain italic text.
$ws.Range("A66").Value = "http://www.ubos.org/onlinefiles/uploads/ubos/pdf%20documents/20067UBR%20report.pdf"
$ws.Range("A66").Font.Italic = $true

# Row 68 no longer holds content (the "UBS" label moves to row 69).
$ws.Range("A68").Clear()

# --- "UBS" citation block (rows 69-72) -----------------------------------
# Each of these rows previously held something formatted the other way
# round (bold/italic alternate every row), so every one needs its font
# explicitly corrected in addition to the new "UBS" text.
$ws.Range("A69").Value = "UBS"
$ws.Range("A69").Font.Bold = $true
$ws.Range("A69").Font.Italic = $false

$ws.Range("A70").Value = "UBS"
$ws.Range("A70").Font.Bold = $false
$ws.Range("A70").Font.Italic = $true

$ws.Range("A71").Value = "UBS"
$ws.Range("A71").Font.Bold = $true
$ws.Range("A71").Font.Italic = $false

$ws.Range("A72").Value = "UBS"
$ws.Range("A72").Font.Italic = $true
